$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

# Clear existing hyperlinks; we will re-add them in row order so relationship IDs line up
$ws.Hyperlinks.Delete()

# Row 2
$ws.Range("A2").Value = '2025-10-15 12:38:42'
$ws.Range("B2").Value = 'n8n×Python×AIで公開レポート自動探索・抽出・分類・登録フロー(PoC開発)'
$ws.Range("C2").Value = 'システム開発'
$ws.Range("D2").Value = '100,000 円 ~ 200,000 円 / 固定'
$ws.Range("E2").Value = '期限情報なし'
$ws.Range("F2").Value = 'https://www.lancers.jp/work/detail/5413825'
$ws.Range("G2").Value = 538
$ws.Range("H2").Value = '🔥AI,Python ◆開発'
$ws.Hyperlinks.Add($ws.Range("F2"), 'https://www.lancers.jp/work/detail/5413825') | Out-Null

# Row 3
$ws.Range("A3").Value = '2025-10-15 12:38:42'
$ws.Range("B3").Value = '【急募】AI×LINE開発をリード!医療機関向けアプリのサーバーサイドエンジニア募集(フルリモート)'
$ws.Range("C3").Value = 'システム開発'
$ws.Range("D3").Value = '500,000 円 ~ 1,000,000 円 / 固定'
$ws.Range("E3").Value = '期限情報なし'
$ws.Range("F3").Value = 'https://www.lancers.jp/work/detail/5413230'
$ws.Range("G3").Value = 385
$ws.Range("H3").Value = '🔥AI,Ai ◆開発 ◇アプリ'
$ws.Hyperlinks.Add($ws.Range("F3"), 'https://www.lancers.jp/work/detail/5413230') | Out-Null

# Row 4
$ws.Range("A4").Value = '2025-10-15 12:38:42'
$ws.Range("B4").Value = '【急募】不動産評価システムのAIチャットポット開発者募集'
$ws.Range("C4").Value = 'システム開発'
$ws.Range("D4").Value = '100,000 円 ~ 200,000 円 / 固定'
$ws.Range("E4").Value = '期限情報なし'
$ws.Range("F4").Value = 'https://www.lancers.jp/work/detail/5413280'
$ws.Range("G4").Value = 383
$ws.Range("H4").Value = '🔥AI,Ai ◆開発'
$ws.Hyperlinks.Add($ws.Range("F4"), 'https://www.lancers.jp/work/detail/5413280') | Out-Null

# Row 5
$ws.Range("A5").Value = '2025-10-15 12:38:42'
$ws.Range("B5").Value = '【急募】ローカルAI開発プロジェクトの協力者を探しています!'
$ws.Range("C5").Value = 'システム開発'
$ws.Range("D5").Value = '300,000 円 ~ 500,000 円 / 固定'
$ws.Range("E5").Value = '期限情報なし'
$ws.Range("F5").Value = 'https://www.lancers.jp/work/detail/5413402'
$ws.Range("G5").Value = 375
$ws.Range("H5").Value = '🔥AI,Ai ◆開発'
$ws.Hyperlinks.Add($ws.Range("F5"), 'https://www.lancers.jp/work/detail/5413402') | Out-Null

# Row 6
$ws.Range("A6").Value = '2025-10-15 12:38:42'
$ws.Range("B6").Value = '【高報酬/リモート可/法人可】グローバルHRベンチャーでAIを活用し業務効率化を推進してくださる方!'
$ws.Range("C6").Value = 'システム開発'
$ws.Range("D6").Value = '5,000 円 ~ 10,000 円 / 固定'
$ws.Range("E6").Value = '期限情報なし'
$ws.Range("F6").Value = 'https://www.lancers.jp/work/detail/5413210'
$ws.Range("G6").Value = 370
$ws.Range("H6").Value = '🔥AI,Ai ◆効率化'
$ws.Hyperlinks.Add($ws.Range("F6"), 'https://www.lancers.jp/work/detail/5413210') | Out-Null

# Row 7
$ws.Range("A7").Value = '2025-10-15 12:38:42'
$ws.Range("B7").Value = '【急募】AIテキスト抜粋アプリのプロンプト最適化依頼'
$ws.Range("C7").Value = 'システム開発'
$ws.Range("D7").Value = '100,000 円 ~ 200,000 円 / 固定'
$ws.Range("E7").Value = '期限情報なし'
$ws.Range("F7").Value = 'https://www.lancers.jp/work/detail/5413215'
$ws.Range("G7").Value = 338
$ws.Range("H7").Value = '🔥AI,Ai ◇アプリ'
$ws.Hyperlinks.Add($ws.Range("F7"), 'https://www.lancers.jp/work/detail/5413215') | Out-Null

# Row 8
$ws.Range("A8").Value = '2025-10-15 12:38:42'
$ws.Range("B8").Value = '3Dプリント用データのWeb自動チェック&変換&カラー補正ツール|開発パートナー募集'
$ws.Range("C8").Value = 'システム開発'
$ws.Range("D8").Value = '5,000,000 円 ~ / 固定'
$ws.Range("E8").Value = '期限情報なし'
$ws.Range("F8").Value = 'https://www.lancers.jp/work/detail/5413508'
$ws.Range("G8").Value = 135
$ws.Range("H8").Value = '◆ツール,開発'
$ws.Hyperlinks.Add($ws.Range("F8"), 'https://www.lancers.jp/work/detail/5413508') | Out-Null

# Row 9
$ws.Range("A9").Value = '2025-10-15 12:38:42'
$ws.Range("B9").Value = 'IB報酬を得るための高性能EA開発依頼'
$ws.Range("C9").Value = 'システム開発'
$ws.Range("D9").Value = '100,000 円 ~ 200,000 円 / 固定'
$ws.Range("E9").Value = '期限情報なし'
$ws.Range("F9").Value = 'https://www.lancers.jp/work/detail/5413293'
$ws.Range("G9").Value = 68
$ws.Range("H9").Value = '◆開発'
$ws.Hyperlinks.Add($ws.Range("F9"), 'https://www.lancers.jp/work/detail/5413293') | Out-Null

# Row 10
$ws.Range("A10").Value = '2025-10-15 12:38:42'
$ws.Range("B10").Value = 'wordpressレンダリングを妨げるリソースの除外'
$ws.Range("C10").Value = 'システム開発'
$ws.Range("D10").Value = '200,000 円 ~ 300,000 円 / 固定'
$ws.Range("E10").Value = '期限情報なし'
$ws.Range("F10").Value = 'https://www.lancers.jp/work/detail/5016989'
$ws.Range("G10").Value = 33
$ws.Range("H10").Value = '○WordPress'
$ws.Hyperlinks.Add($ws.Range("F10"), 'https://www.lancers.jp/work/detail/5016989') | Out-Null

# Row 11
$ws.Range("A11").Value = '2025-10-15 12:38:42'
$ws.Range("B11").Value = '【急募】Cloud RunでWordPress構築のプロを探しています!'
$ws.Range("C11").Value = 'システム開発'
$ws.Range("D11").Value = '10,000 円 ~ 20,000 円 / 固定'
$ws.Range("E11").Value = '期限情報なし'
$ws.Range("F11").Value = 'https://www.lancers.jp/work/detail/5413043'
$ws.Range("G11").Value = 25
$ws.Range("H11").Value = '○WordPress'
$ws.Hyperlinks.Add($ws.Range("F11"), 'https://www.lancers.jp/work/detail/5413043') | Out-Null

# Row 12
$ws.Range("A12").Value = '2025-10-15 12:38:42'
$ws.Range("B12").Value = 'Access 32bitから64bitへの修正改善依頼'
$ws.Range("C12").Value = 'システム開発'
$ws.Range("D12").Value = '5,000 円 ~ 10,000 円 / 固定'
$ws.Range("E12").Value = '期限情報なし'
$ws.Range("F12").Value = 'https://www.lancers.jp/work/detail/5413333'
$ws.Range("G12").Value = 10
$ws.Range("H12").ClearContents()
$ws.Hyperlinks.Add($ws.Range("F12"), 'https://www.lancers.jp/work/detail/5413333') | Out-Null

# Normalize hyperlink cell style back to the workbook's Hyperlink cell style
$ws.Range("F2:F12").Style = "Hyperlink"

Write-Output "done"